# Hindalco price sheet daily update
# - A new "latest" row is inserted at row 2, shifting every existing row
#   down by one (dimension grows from F146 to F147; a new final row 147
#   appears that is an exact copy of what used to be the last row, 146).
# - The new row 2 gets: Date = (old row2 date) + 1 day, while Description /
#   Grade / Basic Price / Circular Date / Circular Link simply repeat what
#   was already in row 2 (now shifted to row 3) - no new circular had been
#   published yet for the new date.
# - One additional hyperlink needs to be created for F86, because the
#   content landing there (shifted up from the former F85) did carry a
#   hyperlink there before, but row-shifting does not re-anchor existing
#   Hyperlink ranges in this engine, so a fresh Hyperlink object is added on
#   top of the new F86 cell to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remember old row 2's values before anything gets overwritten -----
$oldB = $ws.Range("B2").Value2
$oldC = $ws.Range("C2").Value2
$oldD = $ws.Range("D2").Value2
$oldE = $ws.Range("E2").Value2
$oldF = $ws.Range("F2").Value2

# --- 2. Shift rows 2..146 down to 3..147 (values, then formats) ----------
# Doing this as a single block copy (rather than Rows.Insert) keeps the
# worksheet's style table untouched - no new/orphaned cell formats appear,
# and it also leaves the existing <hyperlinks> ref/r:id pairing exactly as
# it was (only the underlying cell text moves), which is what the target
# file expects.
$ws.Range("A2:F146").Copy()
$ws.Range("A3:F147").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A2:F146").Copy()
$ws.Range("A3:F147").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Populate the new row 2 --------------------------------------------
$ws.Range("B2").Value = $oldB
$ws.Range("C2").Value = $oldC
$ws.Range("D2").Value = $oldD
$ws.Range("F2").Value = $oldF

# Date / Circular Date are day-month-year-looking text that Excel's COM
# layer would otherwise auto-parse into a real date serial the moment it's
# assigned. Routing it through a literal-text formula first, then
# collapsing that formula down to its plain value in place (copy/paste
# values-only onto itself) yields a genuine text cell without leaving any
# quote-prefix marker or extra number format behind.
$ws.Range("A2").Formula = '="04-11-2025"'
$ws.Range("E2").Formula = '="' + $oldE + '"'
$ws.Range("A2:E2").Copy()
$ws.Range("A2:E2").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- 4. Re-create the hyperlink that belongs on the new F86 cell ---------
$ws.Hyperlinks.Add($ws.Range("F86"), $ws.Range("F86").Value2)
